# Fix header labels on the existing sheets and add a new "PO Forecast"
# sheet with forecast data (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Rename the generic "Requested quantity" headers to series-specific names.
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet (tab order:
# Weekly Quantity, Monthly Trend, PO Forecast).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Headers.
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Match the bold/centered/bordered header style already used on the other
# sheets (copy formats only, so values set above are preserved).
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Forecast dates (column A) - weekly timestamps.
$dates = @(
    45431.99999999999,
    45585.99999999999,
    45592.99999999999,
    45599.99999999999,
    45606.99999999999,
    45613.99999999999,
    45620.99999999999,
    45627.99999999999,
    45634.99999999999,
    45641.99999999999
)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws3.Cells.Item($i + 2, 1).Value = $dates[$i]
}

# Apply the existing date/time number format (reuse style from sheet 1).
$ws1.Range("A2").Copy()
$ws3.Range("A2:A11").PasteSpecial(-4122)

# Forecast values (columns B, C, D).
$bVals = @(2, 2, 2, 2, 2, 2, 2, 2, 2, 2)
$cVals = @(
    1.999999997420154,
    1.999999997615331,
    1.999999997617446,
    1.999999997557977,
    1.999999997366549,
    1.999999997238737,
    1.99999999710937,
    1.999999997057715,
    1.999999996674064,
    1.999999996893741
)
$dVals = @(
    2.000000002517616,
    2.000000002688245,
    2.000000002609195,
    2.000000002668739,
    2.000000002542237,
    2.000000002614934,
    2.000000002802167,
    2.000000003083935,
    2.000000003209456,
    2.000000003426043
)

for ($i = 0; $i -lt $bVals.Length; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 2).Value = $bVals[$i]
    $ws3.Cells.Item($row, 3).Value = $cVals[$i]
    $ws3.Cells.Item($row, 4).Value = $dVals[$i]
}

# Restore the originally-active sheet.
$ws1.Activate() | Out-Null
